$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 337, shifting existing rows 337:355 down to 338:356.
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new weekly data point.
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = "Femacal de La Calera"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44706
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100112043
$ws.Range("G337").Value = "Pepino ensalada"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 78
$ws.Range("K337").Value = 17000
$ws.Range("L337").Value = 18000
$ws.Range("M337").Value = 17513
$ws.Range("N337").Value = "$/caja 70 unidades"
$ws.Range("O337").Value = "Región de Arica y Parinacota"
$ws.Range("P337").Value = 250
$ws.Range("Q337").Value = 70
$ws.Range("R337").Value = "Hortaliza"
